# Data Refresh (2nd May) — append two new postal-code location rows
# to the "master-location" sheet, mirroring the existing Arabic postal
# code rows (e.g. rows 118/119) for location codes 10113 and 10114.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# code, name (both equal to the numeric location code)
$newCodes = @(10113, 10114)

$startRow = 120
for ($i = 0; $i -lt $newCodes.Length; $i++) {
    $r = $startRow + $i
    $code = $newCodes[$i]

    $ws.Cells.Item($r, 1).Value = $code            # code
    $ws.Cells.Item($r, 2).Value = $code            # name
    $ws.Cells.Item($r, 3).Value = 5                # hierarchy_level
    $ws.Cells.Item($r, 4).Value = "الرمز البريدي"   # hierarchy_level_name (ara - Postal Code)
    $ws.Cells.Item($r, 5).Value = "BNMR"           # parent_loc_code
    $ws.Cells.Item($r, 6).Value = "ara"            # lang_code
    $ws.Cells.Item($r, 7).Value = $true            # is_active
    $ws.Cells.Item($r, 8).Value = "superadmin"     # cr_by
    $ws.Cells.Item($r, 9).Value = "now()"          # cr_dtimes
}

# Reflect the post-edit selection left on the sheet (rows below the
# newly-appended data, as seen in the saved workbook).
[void]$ws.Range("A122:XFD1048576").Select()
